# Update the "talk-schedule" sheet with the new talk's details, switch the
# date cells to plain text ("dd/mm/yyyy") instead of Excel date serials, and
# blank out the rows that held the old (now-missing) schedule entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("talk-schedule")

# --- Row 2: speaker / accommodation / visit details -----------------------
$ws.Range("A2").Value = "Ezzat Elokda"
$ws.Range("B2").Value = "elokdae@ethz.ch"
$ws.Range("C2").Value = "ETH, Zurich"
$ws.Range("D2").Value = "Name: Four Points Flex by Sheraton Lyngby"
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = "18/09/2025"
$ws.Range("G2").Value = "Shobhit Singhal"
$ws.Range("H2").Value = "11:00 - 12:00"
$ws.Range("I2").Value = "18/09/2025"
$ws.Range("J2").Value = "R113"
$ws.Range("K2").Value = "Shobhit Singhal"
$ws.Range("L2").Value = "Room: LY325-R113"

# --- Row 3: accommodation address / seminar host / time -------------------
$ws.Range("D3").Value = "Address: Lundtoftegårdsvej 12, 2800 Kongens Lyngby"
$ws.Range("G3").Value = "Jalal Kazempour"
$ws.Range("H3").Value = "13:00 - 14:00"
$ws.Range("I3").Value = "18/09/2025"
$ws.Range("L3").Value = "Time: 10:00 - 11:00"

# --- Row 4: accommodation map link / seminar date --------------------------
$ws.Range("D4").Value = "Google Maps: https://share.google/o8XML8eIcPmbasHbG"
$ws.Range("L4").Value = "Date: 18/09/2025"

# --- Row 5: seminar title ---------------------------------------------------
$ws.Range("L5").Value = "Title: Aligning the ""Socio"" in Socio-Technical Control: Trustworthy, Fair, and Efficient Resource Allocation with Karma Economies"

# --- Clear the old (now missing) schedule rows in columns G:J --------------
$ws.Range("G4:J10").ClearContents()

$ws.Range("D5").Select() | Out-Null
